$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old row of values (C1:P1) so the sheet dimension shrinks to A1:B1
$ws.Range("C1:P1").Clear()

# Set the new random, non-overlapping reaction values
$ws.Range("A1").Value = 34
$ws.Range("B1").Value = 35
